# Update Name of Algo
# Applies the numeric corrections described by the commit diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "A3"   = -21.986
    "A14"  = -21.759
    "A21"  = -20.269
    "C22"  = -12.711
    "A23"  = -20.246
    "C24"  = -12.218
    "A25"  = -21.855
    "A26"  = -22.066
    "C28"  = -13.179
    "A29"  = -21.311
    "C36"  = -12.732
    "C45"  = -12.921
    "C48"  = -11.205
    "C49"  = -13.367
    "C52"  = -11.655
    "A53"  = -21.836
    "C53"  = -12.789
    "C54"  = -13.061
    "A57"  = -22.17
    "A59"  = -22.406
    "A69"  = -21.519
    "C70"  = -11.492
    "A79"  = -21.008
    "A83"  = -21.982
    "C86"  = -13.9
    "C87"  = -12.907
    "C89"  = -13.299
    "A91"  = -20.744
    "A93"  = -21.508
    "C101" = -12.721
    "A103" = -22.052
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$wb.Save()
